$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 117.066666
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H11").Value = 1009.8
$ws.Range("I11").Value = 1009.8
$ws.Range("K11").Value = 1009.8
$ws.Range("M11").Value = -869.8
$ws.Range("H93").Value = 15000
$ws.Range("J93").Value = 15000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -19992
$ws.Range("H113").Value = 1781.6666
$ws.Range("I113").Value = 1807.2727
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1807.2727
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1446.7273
$ws.Range("N113").Value = -8008
$ws.Range("H132").Value = 7577.5
$ws.Range("I132").Value = 7577.5
$ws.Range("K132").Value = 22732.5
$ws.Range("M132").Value = -20202.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4005.647
$ws.Range("I32").Value = 4267.5713
$ws.Range("J32").Value = 2783.3333
$ws.Range("K32").Value = 4267.5713
$ws.Range("L32").Value = 2783.3333
$ws.Range("M32").Value = -3980.5713
$ws.Range("N32").Value = -3357.3333
$ws.Range("H95").Value = 14425.667
$ws.Range("J95").Value = 14425.667
$ws.Range("L95").Value = 14425.667
$ws.Range("N95").Value = -19917.667
$ws.Range("H96").Value = 30948.166
$ws.Range("J96").Value = 30948.166
$ws.Range("L96").Value = 30948.166
$ws.Range("N96").Value = -36440.166
$ws.Range("H98").Value = 20354.5
$ws.Range("J98").Value = 20354.5
$ws.Range("L98").Value = 20354.5
$ws.Range("N98").Value = -26344.5
$ws.Range("H122").Value = 11597.2
$ws.Range("I122").Value = 11597.2
$ws.Range("K122").Value = 34791.60000000001
$ws.Range("M122").Value = -32341.60000000001
$ws.Range("H132").Value = 1542.4
$ws.Range("I132").Value = 928
$ws.Range("K132").Value = 2784
$ws.Range("M132").Value = -254

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 460.16666
$ws.Range("J64").Value = 457
$ws.Range("L64").Value = 457
$ws.Range("N64").Value = -907
$ws.Range("H67").Value = 460.16666
$ws.Range("J67").Value = 457
$ws.Range("L67").Value = 457
$ws.Range("N67").Value = -2017

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1464
$ws.Range("I6").Value = 284.66666
$ws.Range("K6").Value = 284.66666
$ws.Range("M6").Value = -171.66666
$ws.Range("H12").Value = 918.75
$ws.Range("J12").Value = 918.75
$ws.Range("L12").Value = 918.75
$ws.Range("N12").Value = -1258.75
$ws.Range("H17").Value = 2750
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4826
$ws.Range("H28").Value = 41199.5
$ws.Range("J28").Value = 41199.5
$ws.Range("L28").Value = 41199.5
$ws.Range("N28").Value = -41689.5
$ws.Range("H31").Value = 2200.4
$ws.Range("I31").Value = 2001.1666
$ws.Range("J31").Value = 2499.25
$ws.Range("K31").Value = 2001.1666
$ws.Range("L31").Value = 2499.25
$ws.Range("M31").Value = -1706.1666
$ws.Range("N31").Value = -3089.25
$ws.Range("H34").Value = 2200.4
$ws.Range("I34").Value = 2001.1666
$ws.Range("J34").Value = 2499.25
$ws.Range("K34").Value = 2001.1666
$ws.Range("L34").Value = 2499.25
$ws.Range("M34").Value = -1799.1666
$ws.Range("N34").Value = -2903.25
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
$ws.Range("H99").Value = 33937.5
$ws.Range("I99").Value = 31916.666
$ws.Range("K99").Value = 31916.666
$ws.Range("M99").Value = -30418.666
$ws.Range("H126").Value = 33937.5
$ws.Range("I126").Value = 31916.666
$ws.Range("K126").Value = 95749.99800000001
$ws.Range("M126").Value = -93279.99800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 3173.25
$ws.Range("I25").Value = 224
$ws.Range("J25").Value = 6122.5
$ws.Range("K25").Value = 672
$ws.Range("L25").Value = 18367.5
$ws.Range("M25").Value = -503
$ws.Range("N25").Value = -18705.5
$ws.Range("H30").Value = 3173.25
$ws.Range("I30").Value = 224
$ws.Range("J30").Value = 6122.5
$ws.Range("K30").Value = 672
$ws.Range("L30").Value = 18367.5
$ws.Range("M30").Value = -570
$ws.Range("N30").Value = -18571.5
$ws.Range("H121").Value = 1205
$ws.Range("I121").Value = 1205
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 3615
$ws.Range("L121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -2305
$ws.Range("H124").Value = 5000
$ws.Range("I124").Value = 5000
$ws.Range("K124").Value = 15000
$ws.Range("M124").Value = -10090

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 995
$ws.Range("I70").Value = 995
$ws.Range("K70").Value = 995
$ws.Range("M70").Value = -725
$ws.Range("H73").Value = 995
$ws.Range("I73").Value = 995
$ws.Range("K73").Value = 995
$ws.Range("M73").Value = -59
$ws.Range("H92").Value = 6550
$ws.Range("J92").Value = 6550
$ws.Range("L92").Value = 6550
$ws.Range("N92").Value = -10294
$ws.Range("H101").Value = 19999
$ws.Range("J101").Value = 19999
$ws.Range("L101").Value = 19999
$ws.Range("N101").Value = -26489
$ws.Range("H102").Value = 1558.3334
$ws.Range("I102").Value = 1558.3334
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1558.3334
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 63.66660000000002
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550
$ws.Range("H126").Value = 3001
$ws.Range("I126").Value = 3001
$ws.Range("K126").Value = 9003
$ws.Range("M126").Value = -6533

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H40").Value = 8447.75
$ws.Range("I40").Value = 8630.333000000001
$ws.Range("J40").Value = 7900
$ws.Range("K40").Value = 8630.333000000001
$ws.Range("L40").Value = 7900
$ws.Range("M40").Value = -8494.333000000001
$ws.Range("N40").Value = -8172
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2765.6667
$ws.Range("I122").Value = 2765.6667
$ws.Range("K122").Value = 8297.000100000001
$ws.Range("M122").Value = -5847.000100000001
